$wb = $excel.ActiveWorkbook

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 491.66666
$ws.Range("I19").Value = 189.66667
$ws.Range("K19").Value = 189.66667
$ws.Range("M19").Value = -14.66667000000001

# ALC!row69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3506.5
$ws.Range("I69").Value = 3506.5
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 10519.5
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -9645.5
$ws.Range("N69").ClearContents()

# ALC!row72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 3506.5
$ws.Range("I72").Value = 3506.5
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 31558.5
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -27190.5
$ws.Range("N72").ClearContents()

# ALC!row100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 7851.706
$ws.Range("I100").Value = 12567.223
$ws.Range("J100").Value = 2546.75
$ws.Range("K100").Value = 12567.223
$ws.Range("L100").Value = 2546.75
$ws.Range("M100").Value = -12026.223
$ws.Range("N100").Value = -3628.75

# ALC!row131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2613.6316
$ws.Range("I131").Value = 808
$ws.Range("J131").Value = 4619.8887
$ws.Range("K131").Value = 2424
$ws.Range("L131").Value = 13859.6661
$ws.Range("M131").Value = 2616
$ws.Range("N131").Value = -23939.6661

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 12201790
$ws.Range("I132").Value = 13519109
$ws.Range("J132").Value = 16587.5
$ws.Range("K132").Value = 40557327
$ws.Range("L132").Value = 49762.5
$ws.Range("M132").Value = -40554797
$ws.Range("N132").Value = -54822.5

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25646892
$ws.Range("I32").Value = 14930082
$ws.Range("K32").Value = 14930082
$ws.Range("M32").Value = -14929795

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1949.6842
$ws.Range("I61").Value = 1904.7
$ws.Range("J61").Value = 1999.6666
$ws.Range("K61").Value = 1904.7
$ws.Range("L61").Value = 1999.6666
$ws.Range("M61").Value = -1692.7
$ws.Range("N61").Value = -2423.6666

# ARM!row111
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1949.6842
$ws.Range("I136").Value = 1904.7
$ws.Range("J136").Value = 1999.6666
$ws.Range("K136").Value = 5714.1
$ws.Range("L136").Value = 5998.9998
$ws.Range("M136").Value = -3164.1
$ws.Range("N136").Value = -11098.9998

# BSM!row53
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 48000
$ws.Range("J53").Value = 48000
$ws.Range("L53").Value = 48000
$ws.Range("N53").Value = -49148

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 988.4
$ws.Range("I99").Value = 852.8946999999999
$ws.Range("J99").Value = 1222.4546
$ws.Range("K99").Value = 852.8946999999999
$ws.Range("L99").Value = 1222.4546
$ws.Range("M99").Value = 645.1053000000001
$ws.Range("N99").Value = -4218.4546

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2781627.5
$ws.Range("I134").Value = 1119.2858
$ws.Range("J134").Value = 5854821
$ws.Range("K134").Value = 3357.8574
$ws.Range("L134").Value = 17564463
$ws.Range("M134").Value = -822.8574000000003
$ws.Range("N134").Value = -17569533

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1502.6884
$ws.Range("I31").Value = 1217
$ws.Range("J31").Value = 1632.0566
$ws.Range("K31").Value = 1217
$ws.Range("L31").Value = 1632.0566
$ws.Range("M31").Value = -922
$ws.Range("N31").Value = -2222.0566

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1502.6884
$ws.Range("I34").Value = 1217
$ws.Range("J34").Value = 1632.0566
$ws.Range("K34").Value = 1217
$ws.Range("L34").Value = 1632.0566
$ws.Range("M34").Value = -1015
$ws.Range("N34").Value = -2036.0566

# CUL!row12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 33357.79
$ws.Range("I12").Value = 3.4444444
$ws.Range("J12").Value = 45865.668
$ws.Range("K12").Value = 10.3333332
$ws.Range("L12").Value = 137597.004
$ws.Range("M12").Value = 162.6666668
$ws.Range("N12").Value = -137943.004

# CUL!row45
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 185
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 185
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 555
$ws.Range("N45").Value = -1619
$ws.Range("M45").ClearContents()

# CUL!row50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 2301.5
$ws.Range("J50").Value = 2301.5
$ws.Range("L50").Value = 6904.5
$ws.Range("N50").Value = -7866.5

# CUL!row53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 2301.5
$ws.Range("J53").Value = 2301.5
$ws.Range("L53").Value = 6904.5
$ws.Range("N53").Value = -7866.5

# CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6848.5
$ws.Range("I68").Value = 516.5
$ws.Range("J68").Value = 10647.7
$ws.Range("K68").Value = 1549.5
$ws.Range("L68").Value = 31943.1
$ws.Range("M68").Value = -738.5
$ws.Range("N68").Value = -33565.10000000001

# CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 6848.5
$ws.Range("I71").Value = 516.5
$ws.Range("J71").Value = 10647.7
$ws.Range("K71").Value = 4648.5
$ws.Range("L71").Value = 95829.3
$ws.Range("M71").Value = -592.5
$ws.Range("N71").Value = -103941.3

# CUL!row92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 11384.9
$ws.Range("I92").Value = 482
$ws.Range("J92").Value = 16057.571
$ws.Range("K92").Value = 1446
$ws.Range("L92").Value = 48172.713
$ws.Range("M92").Value = -198
$ws.Range("N92").Value = -50668.713

# CUL!row97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 807.8333
$ws.Range("I97").Value = 631
$ws.Range("J97").Value = 984.6667
$ws.Range("K97").Value = 1893
$ws.Range("L97").Value = 2954.0001
$ws.Range("M97").Value = -1397
$ws.Range("N97").Value = -3946.0001

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 55567516
$ws.Range("I122").Value = 250000320
$ws.Range("J122").Value = 15285.571
$ws.Range("K122").Value = 2250002880
$ws.Range("L122").Value = 137570.139
$ws.Range("M122").Value = -2250000430
$ws.Range("N122").Value = -142470.139

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 23816150
$ws.Range("I132").Value = 666
$ws.Range("J132").Value = 25006924
$ws.Range("K132").Value = 5994
$ws.Range("L132").Value = 225062316
$ws.Range("M132").Value = -3464
$ws.Range("N132").Value = -225067376

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1901.3334
$ws.Range("I7").Value = 1881.6
$ws.Range("K7").Value = 1881.6
$ws.Range("M7").Value = -1769.6

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1950
$ws.Range("I82").Value = 1950
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1950
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1589
$ws.Range("N82").ClearContents()

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1950
$ws.Range("I85").Value = 1950
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1950
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -702
$ws.Range("N85").ClearContents()

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1901.3334
$ws.Range("I126").Value = 1881.6
$ws.Range("K126").Value = 5644.799999999999
$ws.Range("M126").Value = -3174.799999999999

# WVR!row45
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10418
$ws.Range("J45").Value = 11487.667
$ws.Range("L45").Value = 11487.667
$ws.Range("N45").Value = -12469.667

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 31250768
$ws.Range("I81").Value = 41667460
$ws.Range("J81").Value = 695
$ws.Range("K81").Value = 83334920
$ws.Range("L81").Value = 1390
$ws.Range("M81").Value = -83333859
$ws.Range("N81").Value = -3512

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 31250768
$ws.Range("I84").Value = 41667460
$ws.Range("J84").Value = 695
$ws.Range("K84").Value = 416674600
$ws.Range("L84").Value = 6950
$ws.Range("M84").Value = -416669296
$ws.Range("N84").Value = -17558
